# New crime data collected — weekly 68th Precinct CompStat update.
# Bump the report volume/number and the covered week's date range, then
# refresh the crime-complaint figures for rows 15-33 (Rape .. Traffic
# Fatalities) with the newly collected counts/percentages.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header: "Volume 32   Number  11" -> "...  12" and the covered week ---
$ws.Range("A8").Value = "Volume 32   Number  12"
$ws.Range("C9").Value = "Report Covering the Week  3/17/2025  Through  3/23/2025"

# --- Row 15 (Rape) ---
# C15 goes from an empty placeholder ("0") to an actual count -> adopt the
# numeric (#,##0) style used by its neighbours before writing the value.
$ws.Range("C15").Value = 1
$ws.Range("C39").Copy()
$ws.Range("C15").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("F15").Value = 2
$ws.Range("I15").Value = 5
$ws.Range("L15").Value = 150
$ws.Range("M15").Value = 400
$ws.Range("N15").Value = -37.5

# --- Row 16 (Robbery) ---
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 4
$ws.Range("H16").Value = -75
$ws.Range("J16").Value = 22
$ws.Range("K16").Value = -63.636363636363
$ws.Range("L16").Value = -60
$ws.Range("M16").Value = -65.217391304347
$ws.Range("N16").Value = -94.244604316546

# --- Row 17 (Fel. Assault) ---
# C17 goes the other way: from a count back down to the "0" placeholder ->
# adopt the text style used by its neighbours.
$ws.Range("C17").Value = "'0"
$ws.Range("C14").Copy()
$ws.Range("C17").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("D17").Value = 3
$ws.Range("E17").Value = -100
$ws.Range("F17").Value = 7
$ws.Range("H17").Value = -30
$ws.Range("J17").Value = 23
$ws.Range("K17").Value = 4.347826086956
$ws.Range("L17").Value = -33.333333333333
$ws.Range("N17").Value = -60

# --- Row 18 (Burglary) ---
$ws.Range("C18").Value = 3
$ws.Range("F18").Value = 4
$ws.Range("G18").Value = 4
$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 21
$ws.Range("K18").Value = -8.695652173913
$ws.Range("L18").Value = -41.666666666666
$ws.Range("M18").Value = -46.153846153846
$ws.Range("N18").Value = -92.250922509225

# --- Row 19 (Gr. Larceny) ---
$ws.Range("C19").Value = 4
$ws.Range("D19").Value = 7
$ws.Range("E19").Value = -42.857142857142
$ws.Range("F19").Value = 30
$ws.Range("G19").Value = 27
$ws.Range("H19").Value = 11.111111111111
$ws.Range("I19").Value = 63
$ws.Range("J19").Value = 99
$ws.Range("K19").Value = -36.363636363636
$ws.Range("L19").Value = -31.521739130434
$ws.Range("M19").Value = -8.695652173913
$ws.Range("N19").Value = -47.058823529411

# --- Row 20 (G.L.A.) ---
$ws.Range("C20").Value = 1
$ws.Range("E20").Value = 0
$ws.Range("G20").Value = 4
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 10
$ws.Range("J20").Value = 48
$ws.Range("K20").Value = -79.166666666666
$ws.Range("L20").Value = -54.545454545454
$ws.Range("M20").Value = -69.696969696969
$ws.Range("N20").Value = -97.747747747747

# --- Row 21 (TOTAL) ---
$ws.Range("C21").Value = 9
$ws.Range("D21").Value = 12
$ws.Range("E21").Value = -25
$ws.Range("F21").Value = 48
$ws.Range("G21").Value = 49
$ws.Range("H21").Value = -2.040816326530
$ws.Range("I21").Value = 131
$ws.Range("J21").Value = 215
$ws.Range("K21").Value = -39.069767441860
$ws.Range("L21").Value = -37.320574162679
$ws.Range("M21").Value = -28.804347826087
$ws.Range("N21").Value = -87.428023032629

# --- Row 22 (Transit): both G22 and H22 drop back to "no data" placeholders ---
$ws.Range("G22").Value = "'0"
$ws.Range("G14").Copy()
$ws.Range("G22").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("H22").Value = "'***.*"
$ws.Range("H14").Copy()
$ws.Range("H22").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Row 24 (Petit Larceny) ---
$ws.Range("C24").Value = 30
$ws.Range("D24").Value = 23
$ws.Range("E24").Value = 30.434782608695
$ws.Range("F24").Value = 95
$ws.Range("G24").Value = 68
$ws.Range("H24").Value = 39.705882352941
$ws.Range("I24").Value = 294
$ws.Range("J24").Value = 292
$ws.Range("K24").Value = 0.684931506849
$ws.Range("L24").Value = -21.6
$ws.Range("M24").Value = 16.205533596837

# --- Row 25 (Retail Theft) ---
$ws.Range("C25").Value = 11
$ws.Range("D25").Value = 17
$ws.Range("E25").Value = -35.294117647058
$ws.Range("F25").Value = 34
$ws.Range("G25").Value = 38
$ws.Range("H25").Value = -10.526315789473
$ws.Range("I25").Value = 146
$ws.Range("J25").Value = 182
$ws.Range("K25").Value = -19.780219780219
$ws.Range("L25").Value = -27

# --- Row 26 (Misd. Assault) ---
$ws.Range("C26").Value = 12
$ws.Range("D26").Value = 6
$ws.Range("E26").Value = 100
$ws.Range("F26").Value = 30
$ws.Range("G26").Value = 25
$ws.Range("H26").Value = 20
$ws.Range("I26").Value = 86
$ws.Range("J26").Value = 72
$ws.Range("K26").Value = 19.444444444444
$ws.Range("L26").Value = -1.149425287356
$ws.Range("M26").Value = 8.860759493670

# --- Row 27 (UCR Rape*) ---
# C27 goes from the "0" placeholder to an actual numeric count.
$ws.Range("C27").Value = 2
$ws.Range("C39").Copy()
$ws.Range("C27").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("F27").Value = 4
$ws.Range("H27").Value = 300
$ws.Range("I27").Value = 7
$ws.Range("K27").Value = 600
$ws.Range("L27").Value = 133.333333333333

# --- Row 28 (Other Sex Crimes) ---
$ws.Range("G28").Value = 5
$ws.Range("H28").Value = -80

# --- Row 31 (Shooting Inc.) ---
# D31/E31 go from placeholders to real numbers.
$ws.Range("D31").Value = 1
$ws.Range("D16").Copy()
$ws.Range("D31").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("E31").Value = -100
$ws.Range("E16").Copy()
$ws.Range("E31").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("G31").Value = 2
$ws.Range("J31").Value = 2

# --- Row 33 (Traffic Fatalities): D33/E33 drop back to placeholders ---
$ws.Range("D33").Value = "'0"
$ws.Range("D14").Copy()
$ws.Range("D33").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("E33").Value = "'***.*"
$ws.Range("E14").Copy()
$ws.Range("E33").PasteSpecial(-4122)
$excel.CutCopyMode = $false
